$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 389.93332
$ws.Range("I11").Value = 389.93332
$ws.Range("K11").Value = 389.93332
$ws.Range("M11").Value = -249.93332
$ws.Range("H32").Value = 15221.723
$ws.Range("I32").Value = 13887.333
$ws.Range("J32").Value = 16556.111
$ws.Range("K32").Value = 13887.333
$ws.Range("L32").Value = 16556.111
$ws.Range("M32").Value = -13561.333
$ws.Range("N32").Value = -17208.111
$ws.Range("H33").Value = 1484.7273
$ws.Range("I33").Value = 1146.75
$ws.Range("J33").Value = 1677.8572
$ws.Range("K33").Value = 1146.75
$ws.Range("L33").Value = 1677.8572
$ws.Range("M33").Value = -917.75
$ws.Range("N33").Value = -2135.8572
$ws.Range("H42").Value = 761.8
$ws.Range("I42").Value = 827.25
$ws.Range("K42").Value = 2481.75
$ws.Range("M42").Value = -2251.75
$ws.Range("H70").Value = 8240.333000000001
$ws.Range("I70").Value = 14149.375
$ws.Range("J70").Value = 5752.316
$ws.Range("K70").Value = 42448.125
$ws.Range("L70").Value = 17256.948
$ws.Range("M70").Value = -42178.125
$ws.Range("N70").Value = -17796.948
$ws.Range("H73").Value = 8240.333000000001
$ws.Range("I73").Value = 14149.375
$ws.Range("J73").Value = 5752.316
$ws.Range("K73").Value = 42448.125
$ws.Range("L73").Value = 17256.948
$ws.Range("M73").Value = -41512.125
$ws.Range("N73").Value = -19128.948
$ws.Range("H132").Value = 3419519.5
$ws.Range("I132").Value = 3502855.2
$ws.Range("K132").Value = 10508565.6
$ws.Range("M132").Value = -10506035.6
$ws.Range("H137").Value = 8043.231
$ws.Range("I137").Value = 10395.912
$ws.Range("K137").Value = 31187.736
$ws.Range("M137").Value = -28637.736
$ws.Range("H138").Value = 2772.4211
$ws.Range("I138").Value = 1930.091
$ws.Range("K138").Value = 5790.272999999999
$ws.Range("M138").Value = -650.2729999999992
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 701
$ws.Range("I5").Value = 701
$ws.Range("K5").Value = 701
$ws.Range("M5").Value = -589
$ws.Range("H32").Value = 18733.928
$ws.Range("I32").Value = 18733.928
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 18733.928
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -18446.928
$ws.Range("N32").ClearContents()
$ws.Range("H45").Value = 3733.7222
$ws.Range("I45").Value = 3041.5
$ws.Range("J45").Value = 4425.9443
$ws.Range("K45").Value = 3041.5
$ws.Range("L45").Value = 4425.9443
$ws.Range("M45").Value = -2664.5
$ws.Range("N45").Value = -5179.9443
$ws.Range("H102").Value = 6378.5625
$ws.Range("I102").Value = 5753.9165
$ws.Range("K102").Value = 5753.9165
$ws.Range("M102").Value = -4131.9165
$ws.Range("H132").Value = 3355.95
$ws.Range("I132").Value = 2032.6666
$ws.Range("J132").Value = 3923.0715
$ws.Range("K132").Value = 6097.9998
$ws.Range("L132").Value = 11769.2145
$ws.Range("N132").Value = -16829.2145
$ws.Range("M132").Value = -3567.9998
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 701
$ws.Range("I4").Value = 701
$ws.Range("K4").Value = 701
$ws.Range("M4").Value = -586
$ws.Range("H75").Value = 27740.25
$ws.Range("I75").Value = 5999.5
$ws.Range("K75").Value = 5999.5
$ws.Range("M75").Value = -5063.5
$ws.Range("H78").Value = 27740.25
$ws.Range("I78").Value = 5999.5
$ws.Range("K78").Value = 17998.5
$ws.Range("M78").Value = -13318.5
$ws.Range("H86").Value = 1633.75
$ws.Range("I86").Value = 1710.2727
$ws.Range("K86").Value = 1710.2727
$ws.Range("M86").Value = -587.2727
$ws.Range("H89").Value = 1633.75
$ws.Range("I89").Value = 1710.2727
$ws.Range("K89").Value = 8551.363499999999
$ws.Range("M89").Value = -2935.363499999999
$ws.Range("H135").Value = 89408.89
$ws.Range("J135").Value = 89408.89
$ws.Range("L135").Value = 89408.89
$ws.Range("N135").Value = -99548.89
$ws.Range("H137").Value = 200000
$ws.Range("J137").Value = 200000
$ws.Range("L137").Value = 200000
$ws.Range("N137").Value = -210200
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 36995
$ws.Range("I26").Value = 44000
$ws.Range("J26").Value = 29990
$ws.Range("K26").Value = 44000
$ws.Range("L26").Value = 29990
$ws.Range("N26").Value = -30564
$ws.Range("M26").Value = -43713
$ws.Range("H62").Value = 5965.5
$ws.Range("I62").Value = 3219.75
$ws.Range("J62").Value = 7534.5
$ws.Range("K62").Value = 3219.75
$ws.Range("L62").Value = 7534.5
$ws.Range("M62").Value = -2595.75
$ws.Range("N62").Value = -8782.5
$ws.Range("H65").Value = 5965.5
$ws.Range("I65").Value = 3219.75
$ws.Range("J65").Value = 7534.5
$ws.Range("K65").Value = 16098.75
$ws.Range("L65").Value = 37672.5
$ws.Range("M65").Value = -12978.75
$ws.Range("N65").Value = -43912.5
$ws.Range("H105").Value = 6425
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H134").Value = 2488.3076
$ws.Range("I134").Value = 1754.9524
$ws.Range("K134").Value = 5264.857199999999
$ws.Range("M134").Value = -2729.857199999999
$ws.Range("H135").Value = 111899.586
$ws.Range("J135").Value = 111899.586
$ws.Range("L135").Value = 111899.586
$ws.Range("N135").Value = -122039.586
$ws.Range("H137").Value = 85267.664
$ws.Range("J137").Value = 85267.664
$ws.Range("L137").Value = 85267.664
$ws.Range("N137").Value = -95467.664
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1381.4
$ws.Range("I107").Value = 1731
$ws.Range("J107").Value = 1095.3636
$ws.Range("K107").Value = 5193
$ws.Range("L107").Value = 3286.0908
$ws.Range("M107").Value = -3273
$ws.Range("N107").Value = -7126.0908
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 5025722
$ws.Range("J44").Value = 6024200
$ws.Range("L44").Value = 6024200
$ws.Range("N44").Value = -6025392
$ws.Range("H80").Value = 5331.593
$ws.Range("I80").Value = 3332.85
$ws.Range("J80").Value = 11042.286
$ws.Range("K80").Value = 3332.85
$ws.Range("L80").Value = 11042.286
$ws.Range("M80").Value = -2334.85
$ws.Range("N80").Value = -13038.286
$ws.Range("H83").Value = 5331.593
$ws.Range("I83").Value = 3332.85
$ws.Range("J83").Value = 11042.286
$ws.Range("K83").Value = 16664.25
$ws.Range("L83").Value = 55211.43
$ws.Range("M83").Value = -11672.25
$ws.Range("N83").Value = -65195.43
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 60000
$ws.Range("J23").Value = 60000
$ws.Range("L23").Value = 60000
$ws.Range("N23").Value = -60460
$ws.Range("H100").Value = 2889.6191
$ws.Range("I100").Value = 2765.6667
$ws.Range("K100").Value = 2765.6667
$ws.Range("M100").Value = -2224.6667
$ws.Range("H132").Value = 4820.4287
$ws.Range("I132").Value = 3874
$ws.Range("K132").Value = 11622
$ws.Range("M132").Value = -9092
$ws.Range("H136").Value = 2672.7
$ws.Range("I136").Value = 2136.3333
$ws.Range("K136").Value = 6408.999899999999
$ws.Range("M136").Value = -3858.999899999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 94661.5
$ws.Range("J70").Value = 95596.8
$ws.Range("L70").Value = 95596.8
$ws.Range("N70").Value = -96226.8
$ws.Range("H73").Value = 94661.5
$ws.Range("J73").Value = 95596.8
$ws.Range("L73").Value = 95596.8
$ws.Range("N73").Value = -97780.8
$ws.Range("H96").Value = 1775.9333
$ws.Range("I96").Value = 1524.6666
$ws.Range("J96").Value = 1943.4445
$ws.Range("K96").Value = 1524.6666
$ws.Range("L96").Value = 1943.4445
$ws.Range("M96").Value = -151.6666
$ws.Range("N96").Value = -4689.4445
$ws.Range("H100").Value = 1300
$ws.Range("I100").Value = 200
$ws.Range("J100").Value = 1666.6666
$ws.Range("K100").Value = 400
$ws.Range("L100").Value = 3333.3332
$ws.Range("M100").Value = 141
$ws.Range("N100").Value = -4415.3332
$ws.Range("H132").Value = 4160.483
$ws.Range("I132").Value = 4269.1665
$ws.Range("J132").Value = 3638.8
$ws.Range("K132").Value = 12807.4995
$ws.Range("L132").Value = 10916.4
$ws.Range("M132").Value = -10277.4995
$ws.Range("N132").Value = -15976.4
$ws.Range("H136").Value = 19063.475
$ws.Range("I136").Value = 22706.719
$ws.Range("K136").Value = 68120.15700000001
$ws.Range("M136").Value = -65570.15700000001
